$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the row containing the variable "ca_hstime_mnspid" in column B and delete it.
$found = $ws.Range("B:B").Find("ca_hstime_mnspid")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
